# se agrega brecha enge y sin enge vs eeuu
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (old D "brecha_productividad_eeuu_arg" shifts to E,
# old E "brecha_enge_eeuu" shifts to F).
$ws.Columns.Item(4).Insert()

# Insert a new column after the (new) F -> becomes G, for "brecha_sin_enge_eeuu".
$ws.Columns.Item(7).Insert()

# Header row
$ws.Range("D1").Value = "brecha_productividad_sin_enge_tot_level"
$ws.Range("G1").Value = "brecha_sin_enge_eeuu"

# New column D values (brecha_productividad_sin_enge_tot_level), rows 2-29
$colD = @{
    2  = 0.7217866279342013
    3  = 0.7061347470885972
    4  = 0.7205939142286747
    5  = 0.7274295185965005
    6  = 0.7158631339776839
    7  = 0.6964979768282287
    8  = 0.7044691905273722
    9  = 0.6572341654741737
    10 = 0.640122874364606
    11 = 0.4664758337295752
    12 = 0.5188359339706339
    13 = 0.4876241351803389
    14 = 0.4848687707668133
    15 = 0.4933955560261969
    16 = 0.4607912212026656
    17 = 0.4512128566514683
    18 = 0.4811258478933602
    19 = 0.5071842261070042
    20 = 0.5093578642283374
    21 = 0.473622428045198
    22 = 0.4723300233404075
    23 = 0.4421653496897042
    24 = 0.4454741915504722
    25 = 0.4774814540996213
    26 = 0.4121049083080825
    27 = 0.3820579528613384
    28 = 0.4268772236080584
    29 = 0.5274793916539007
}

foreach ($row in $colD.Keys) {
    $ws.Range("D$row").Value = $colD[$row]
}

# New column G values (brecha_sin_enge_eeuu), rows 2-26 only (27-29 stay blank)
$colG = @{
    2  = 0.143635538958906
    3  = 0.1546435096124028
    4  = 0.1549276915591651
    5  = 0.1658539302400021
    6  = 0.1625009314129342
    7  = 0.1455680771570998
    8  = 0.1338491462002007
    9  = 0.1189593839508254
    10 = 0.10562027427016
    11 = 0.06297423755349266
    12 = 0.06692983548221178
    13 = 0.06582925824934575
    14 = 0.0620632026581521
    15 = 0.0646348178394318
    16 = 0.06266760608356252
    17 = 0.0613649485045997
    18 = 0.05917847929088331
    19 = 0.06897705475055257
    20 = 0.06978202739928221
    21 = 0.06015004836174015
    22 = 0.06045824298757217
    23 = 0.05350200731245421
    24 = 0.05702069651846044
    25 = 0.05634281158375531
    26 = 0.04739206445542949
}

foreach ($row in $colG.Keys) {
    $ws.Range("G$row").Value = $colG[$row]
}
# Rows 27-29 keep column G blank (no value set), matching the source data
# which has no brecha_sin_enge_eeuu figure for those years.
